$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.280.79'
$ws.Range('E2').Value = '  -4.79%  '
$ws.Range('D3').Value = '2.242.54'
$ws.Range('E3').Value = '  -5.73%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '320.85'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '101.31'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -7.13%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.585'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -8.08%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.565'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -8.31%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '37.26'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -9.18%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '54.60'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.01%  '
$ws.Range('E12').Value = '  -9.70%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.74'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -9.66%  '
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.586.06'
$ws.Range('E15').Value = '  -5.67%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.869'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -11.84%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.50'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -6.39%  '
$ws.Range('D18').Value = '2.241.61'
$ws.Range('E18').Value = '  -5.60%  '
$ws.Range('D19').Value = '43.213.89'
$ws.Range('E19').Value = '  -4.82%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.47'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -7.89%  '
$ws.Range('E21').Value = '  -8.78%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.56'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -10.66%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.67'
$ws.Range('D23').Style = "Normal"
$ws.Range('E24').Value = '  -14.13%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '238.60'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -8.71%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.19'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -7.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.08'
$ws.Range('D29').Style = "Normal"
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.39'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -16.33%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '36.32'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.39%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0884'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -7.38%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '20.47'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -8.94%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '153.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -8.82%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.72'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -6.29%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.24'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +8.30%  '
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('E39').Value = '  -7.66%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.46'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.14%  '
$ws.Range('E41').Value = '  -10.58%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.69'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -7.99%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0326'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -8.48%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.77'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +6.45%  '
$ws.Range('D46').Value = '1.757.31'
$ws.Range('E46').Value = '  -5.21%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '86.74'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -11.00%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.206'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -10.27%  '
$ws.Range('E49').Value = '  -10.52%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '75.88'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -9.30%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '59.13'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -16.42%  '
